# Version 09-25 - Large changes to qoq and ifoCAST evaluation
# Updates the "revision" column (column B) values on the active sheet
# to reflect the refreshed evaluation-series numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = [double]"-4.440892098500626E-14"
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = [double]"4.440892098500626E-14"
$ws.Range("B6").Value = 0.8761527144223624
$ws.Range("B7").Value = 0.2498186593231866
$ws.Range("B8").Value = -0.2073814777540428
$ws.Range("B9").Value = 0.3094428711141628
$ws.Range("B10").Value = 0.5223161956339206
$ws.Range("B11").Value = 0.3117781822009169
$ws.Range("B12").Value = 0.3783498544355668
$ws.Range("B13").Value = 0.6542703491021484
$ws.Range("B14").Value = -0.6157632300240357
$ws.Range("B15").Value = 0.01566947406670405
$ws.Range("B16").Value = -1.024794128387363
$ws.Range("B17").Value = 0.4283030634637974
$ws.Range("B18").Value = 0.4210175484930634
$ws.Range("B19").Value = 0.3309484344824476
$ws.Range("B20").Value = 0.03534272045342401
$ws.Range("B21").Value = -1.144790596790379
$ws.Range("B22").Value = -0.3737413844400406
$ws.Range("B23").Value = 0.3599799282585359
$ws.Range("B24").Value = 0.5766501347738604
$ws.Range("B25").Value = -0.5004758034602208
$ws.Range("B26").Value = -0.7750463390001627
$ws.Range("B27").Value = 0.2146387198177946
$ws.Range("B28").Value = 0.02351224639369764
$ws.Range("B29").Value = -0.5673139978460418
$ws.Range("B30").Value = 0.03242225351332007
$ws.Range("B31").Value = -0.3625064515619281
$ws.Range("B32").Value = -0.5809412570459083
$ws.Range("B33").Value = 0.3008481533817164
$ws.Range("B34").Value = -0.4092135626299287
$ws.Range("B35").Value = -0.7785741838770672
$ws.Range("B36").Value = -0.9959164522061803
$ws.Range("B37").Value = -0.02772281113909703
$ws.Range("B38").Value = 0.565922940159902
$ws.Range("B39").Value = 0.266820966430148
